$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "608.70") that must remain
# plain text, matching the source data which uses inline/shared strings, not numbers.
# Temporarily force Text format while assigning the value, then clear the format
# again so no stray style index is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "97.567.36"
$ws.Range("E2").Value = "  +4.45%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.125.21"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
Set-TextValue $ws.Range("D5") "241.22"
$ws.Range("E5").Value = "  +1.38%  "

# Row 6
Set-TextValue $ws.Range("D6") "608.70"
$ws.Range("E6").Value = "  -1.08%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.12"
$ws.Range("E7").Value = "  +1.33%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.381"
$ws.Range("E8").Value = "  -2.94%  "

# Row 9
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
Set-TextValue $ws.Range("D10") "3.121.21"
$ws.Range("E10").Value = "  +0.15%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.774"
$ws.Range("E11").Value = "  -7.79%  "

# Row 12
$ws.Range("E12").Value = "  -0.37%  "

# Row 13
Set-TextValue $ws.Range("D13") "97.135.43"
$ws.Range("E13").Value = "  +4.28%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.0000237"
$ws.Range("E14").Value = "  -3.34%  "

# Row 15
Set-TextValue $ws.Range("D15") "5.45"
$ws.Range("E15").Value = "  +0.03%  "

# Row 16
Set-TextValue $ws.Range("D16") "33.54"
$ws.Range("E16").Value = "  -5.01%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.707.50"
$ws.Range("E17").Value = "  +0.28%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.121.03"
$ws.Range("E18").Value = "  -0.01%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D19") "513.08"
$ws.Range("E19").Value = "  +15.71%  "

# Row 20
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D20") "3.37"
$ws.Range("E20").Value = "  -10.51%  "

# Row 21
Set-TextValue $ws.Range("D21") "14.34"
$ws.Range("E21").Value = "  -3.58%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.58"
$ws.Range("E22").Value = "  -7.73%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.0000190"
$ws.Range("E23").Value = "  -5.63%  "

# Row 24
Set-TextValue $ws.Range("D24") "8.66"
$ws.Range("E24").Value = "  -4.88%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "88.09"
$ws.Range("E25").Value = "  +2.63%  "

# Row 26
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D26") "5.42"
$ws.Range("E26").Value = "  -6.29%  "

# Row 27
Set-TextValue $ws.Range("D27") "11.43"
$ws.Range("E27").Value = "  -12.03%  "

# Row 28
Set-TextValue $ws.Range("D28") "3.292.78"
$ws.Range("E28").Value = "  +0.26%  "

# Row 29
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.240"
$ws.Range("E30").Value = "  +0.45%  "

# Row 31
$ws.Range("E31").Value = "  -4.36%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.122"
$ws.Range("E32").Value = "  -2.70%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.996"
$ws.Range("E33").Value = "  -2.81%  "

# Row 34
Set-TextValue $ws.Range("D34") "8.84"
$ws.Range("E34").Value = "  -4.93%  "

# Row 35
Set-TextValue $ws.Range("D35") "26.50"
$ws.Range("E35").Value = "  +1.96%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.150"
$ws.Range("E36").Value = "  -6.96%  "

# Row 37
Set-TextValue $ws.Range("D37") "7.16"
$ws.Range("E37").Value = "  -10.72%  "

# Row 38
Set-TextValue $ws.Range("D38") "1.87"
$ws.Range("E38").Value = "  -2.43%  "

# Row 39
$ws.Range("E39").Value = "  +0.89%  "

# Row 40
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D40") "463.81"
$ws.Range("E40").Value = "  -3.43%  "

# Row 41
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Range("D41") "0.428"
$ws.Range("E41").Value = "  -4.55%  "

# Row 42
Set-TextValue $ws.Range("D42") "1.20"
$ws.Range("E42").Value = "  -7.68%  "

# Row 43
Set-TextValue $ws.Range("D43") "3.51"
$ws.Range("E43").Value = "  -10.01%  "

# Row 44
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
Set-TextValue $ws.Range("D45") "3.05"
$ws.Range("E45").Value = "  -8.28%  "

# Row 46
Set-TextValue $ws.Range("D46") "162.52"
$ws.Range("E46").Value = "  +1.90%  "

# Row 47
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D47") "0.692"
$ws.Range("E47").Value = "  -1.40%  "

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D48") "1.89"
$ws.Range("E48").Value = "  +0.23%  "

# Row 49
Set-TextValue $ws.Range("D49") "4.42"
$ws.Range("E49").Value = "  -0.21%  "

# Row 50
Set-TextValue $ws.Range("D50") "44.11"
$ws.Range("E50").Value = "  +0.22%  "

# Row 51
$ws.Range("E51").Value = "  -0.02%  "
